# Apply the edit described by the diff:
#  - Insert a new row (Cheltenham / Woolworths Southland, 1239 Nepean Highway /
#    31/12/20 6:00pm - 6:30pm / Case shopped at store) after the existing
#    Cheltenham rows, i.e. at sheet row 17 (pushing everything from old row 17
#    down by one).
#  - Remove the Hallam row, the first Moorabbin row (10:45am-12:15pm), the
#    Mordialloc row, and the trailing Wonthaggi row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank row at row 17 (old row 17 "Chirnside Park" and below
#    shift down to row 18+).
$ws.Rows.Item(17).Insert()

# 2. Populate the newly inserted row 17 with the new Cheltenham entry.
$ws.Range("A17").Value = "Cheltenham"
$ws.Range("B17").Value = "Woolworths Southland, 1239 Nepean Highway"
$ws.Range("C17").Value = "31/12/20 6:00pm - 6:30pm"
$ws.Range("D17").Value = "Case shopped at store"

# 3. Delete the rows that were removed, working from the bottom up so row
#    numbers for earlier deletions stay valid.
#    Original rows (before the insert above) were: Hallam=25, Moorabbin
#    (10:45am-12:15pm)=34, Mordialloc=36, Wonthaggi=43. After the insert at
#    row 17 everything from there down shifted by +1, so they are now at
#    26, 35, 37, and 44 respectively.
$ws.Rows.Item(44).Delete()
$ws.Rows.Item(37).Delete()
$ws.Rows.Item(35).Delete()
$ws.Rows.Item(26).Delete()
